$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: new logbook entry (date, duration)
$ws.Range("A13").Value = 43525
$ws.Range("A13").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("B13").Value = "4 uur"

# New "interessante links" entries in column Q
$ws.Range("Q15").Value = "https://console.firebase.google.com/"
$ws.Range("Q16").Value = "https://www.androidhive.info/2015/09/android-material-design-working-with-tabs/"

$ws.Range("Q17").Value = "https://www.androidhive.info/2016/05/android-working-with-card-view-and-recycler-view/"
[void]$ws.Hyperlinks.Add($ws.Range("Q17"), "https://www.androidhive.info/2016/05/android-working-with-card-view-and-recycler-view/")
$ws.Range("Q17").Style = $ws.Range("Q3").Style

# Row 13 description (added after the link cells so shared strings are appended
# in the same order as the authored workbook)
$ws.Range("C13").Value = "Firebase inloggen en wachtwoord vergeten toepassen, viewpager aanmaken, cardview in recyclerview zetten, verder uitwerken schermen, cardview toevoegen"

# Update current selection to match the author's last position
[void]$ws.Range("A14").Select()
